$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (F1) onto the new
# header cell (G1), then set its text. This reuses the existing header
# style (bold font + border) rather than minting a new style entry.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Status"

# Match the author's selection state after the edit.
$ws.Range("G2").Select()
